# Update wiki data pipeline:
# Insert a new linked-page entry "Wayback Machine" immediately above the
# existing "Wikidata" row (row 444), shifting all subsequent rows down by
# one. This grows the used range from A1:A451 to A1:A452.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 444 (pushes "Wikidata" and everything
# below it down one row).
$ws.Rows.Item(444).Insert()

# Populate the newly inserted row with the new page title.
$ws.Range("A444").Value = "Wayback Machine"
